$d = $word.ActiveDocument

# --- Block 1: "Consultar veículos disponíveis" -> "Realizar locação",
#              then add two new bullet paragraphs after it. ---

# 1. Replace the text of the existing paragraph.
$rng = $d.Content
$rng.Find.Execute("- Consultar veículos disponíveis", $true, $false, $false, $false, $false, $true, 1, $false, "- Realizar locação", 2)

# 2. Split a new (empty) paragraph right after "- Realizar locação".
$rng = $d.Content
$rng.Find.Execute("- Realizar locação", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Paragraphs(1).Range.InsertParagraphAfter()

# 3. Fill that new paragraph with "- Consultar clientes" and split another
#    empty paragraph after it.
$rng = $d.Content
$rng.Find.Execute("- Realizar locação", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newPara = $rng.Paragraphs(1).Next().Range
$newPara.Text = "- Consultar clientes"
$newPara.InsertParagraphAfter()

# 4. Fill the following new paragraph with "- Consultar veículos ".
$rng = $d.Content
$rng.Find.Execute("- Consultar clientes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newPara2 = $rng.Paragraphs(1).Next().Range
$newPara2.Text = "- Consultar veículos "

# --- Block 2: "Cliente(..." -> "Cliente (..." (space added),
#              then add a new "Locação (...)" paragraph after it. ---

# 5. Replace the text of the existing paragraph.
$rng = $d.Content
$rng.Find.Execute("- Cliente(Nome, CPF, Sexo, dataNascimento, RG, Nacionalidade, Telefone)", $true, $false, $false, $false, $false, $true, 1, $false, "- Cliente (Nome, CPF, Sexo, dataNascimento, RG, Nacionalidade, Telefone)", 2)

# 6. Split a new (empty) paragraph right after it.
$rng = $d.Content
$rng.Find.Execute("- Cliente (Nome, CPF, Sexo, dataNascimento, RG, Nacionalidade, Telefone)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Paragraphs(1).Range.InsertParagraphAfter()

# 7. Fill the new paragraph with the "Locação" requirement line.
$rng = $d.Content
$rng.Find.Execute("- Cliente (Nome, CPF, Sexo, dataNascimento, RG, Nacionalidade, Telefone)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newPara3 = $rng.Paragraphs(1).Next().Range
$newPara3.Text = "- Locação (dataLocacao, horarioLocacao, dataDevolucao, horarioDevolucao)"
